$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update "配比" (ratio) values in column N, and their mirrored text labels in column T ---
# Each pair: row number, new ratio value, and the exact mirrored label text for column T

$updates = @(
    @{ Row = 6;  Value = 4.174266;  Label = "高返-13_1=[4.174266]" },
    @{ Row = 15; Value = 4.174266;  Label = "高返-13_3=[4.174266]" },
    @{ Row = 28; Value = 3.718342;  Label = "高返-13_4=[3.718342]" },
    @{ Row = 34; Value = 0.5526916; Label = "水洗铁-13_4=[0.5526916]" },
    @{ Row = 39; Value = 0.1988863; Label = "高返-13_5=[0.1988863]" },
    @{ Row = 46; Value = 0.5526916; Label = "水洗铁-13_5=[0.5526916]" },
    @{ Row = 50; Value = 0.2570375; Label = "高返-13_12=[0.2570375]" },
    @{ Row = 57; Value = 1.456512;  Label = "水洗铁-13_12=[1.456512]" },
    @{ Row = 65; Value = 2.076149;  Label = "氧化铁-14-1=[2.076149]" },
    @{ Row = 75; Value = 1.312939;  Label = "氧化铁-14-3=[1.312939]" },
    @{ Row = 89; Value = 0;         Label = "巴西粗粉-14-4=[0.0]" },
    @{ Row = 90; Value = 15;        Label = "高品澳粉-14-12=[15.0]" },
    @{ Row = 93; Value = 14.84864;  Label = "塞拉利昂粉-14-12=[14.84864]" }
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Cells.Item($r, 14).Value = $u.Value    # column N = 14 ("配比")
    $ws.Cells.Item($r, 20).Value = $u.Label    # column T = 20 (mirrored text label)
}

# --- Downstream recalculated summary cells (no formulas in sheet; values stored directly) ---
$ws.Range("M99").Value = 489.2533525189781
$ws.Range("M100").Value = 489.2533525189781
$ws.Range("M101").Value = 489.2533525189781

$ws.Range("C102").Value = 55.90000090477773
$ws.Range("D102").Value = 4.5000000675375
$ws.Range("E102").Value = 2.897264995189999
$ws.Range("F102").Value = 0.9762019118428573
$ws.Range("G102").Value = 3.00000014589
$ws.Range("H102").Value = 0.1632088484
$ws.Range("I102").Value = 0.07965211428500001
$ws.Range("J102").Value = 0.16721795398
